$d = $word.ActiveDocument

# Anchor on the current first paragraph ("C:\>cd NewProject") and open up
# two new empty paragraphs right before it.
$anchor = $d.Paragraphs.Item(1).Range
$anchor.InsertParagraphBefore()
$anchor.InsertParagraphBefore()

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# First new paragraph: "List of commands for github:"
$para1Xml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t xml:space="preserve">List of commands for </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>github</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs.Item(1).Range.InsertXML($para1Xml)

# Second new paragraph: "Create folder Newproject on C drive"
$para2Xml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t xml:space="preserve">Create folder </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Newproject</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> on C drive</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs.Item(2).Range.InsertXML($para2Xml)

Write-Host "Inserted two introductory paragraphs before the command list."
